$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample data row (row 2) entirely.
$ws.Rows.Item(2).Delete()

# Remove the now-orphaned hyperlink definition.
foreach ($h in $ws.Hyperlinks) {
  $h.Delete()
}

# Trim trailing spaces from a few header labels.
$ws.Range("A1").Value = "Employee Code"
$ws.Range("F1").Value = "Designation"
$ws.Range("U1").Value = "ESIC Employee"

# Reset the view: select B7 (clears the old topLeftCell/activeCell state).
$ws.Range("B7").Select() | Out-Null
